$d = $word.ActiveDocument

# Replace the "Denumire: ..." paragraph text with "Persoana fizica"
$d.Content.Find.Execute("Denumire: ESX INTEL WORLD S.R.L.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Persoana fizica", 2)

# Remove the following four paragraphs: CUI/Tax ID, Adresa/Adress, Registrul comertului, Email
$targets = @(
    "CUI/Tax ID no: 38151434",
    "Adresa/Adress: MUNICIPIUL BUCUREŞTI, SECTOR 1, CALEA FLOREASCA, NR.169, CORP X, PARTER, CAMERA 24A",
    "Registrul comertului/Registration no: J40/10130/2019",
    "Email: gyds@dsjb"
)

foreach ($t in $targets) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r", "`a") -eq $t) {
            $p.Range.Delete()
            break
        }
    }
}
